$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows (Balance Sheet section) whose D:J quarter columns are cleared to 0
# (K column retains its existing "NA" placeholder / value, unchanged)
$zeroRows = 41,43,44,45,46,47,48,49,52,54,57,58,59,60,62,66,72,76
foreach ($r in $zeroRows) {
    $addr = "D" + $r + ":J" + $r
    $ws.Range($addr).Value = 0
}

# Row 61 (Long Term Debt) only has its D column cleared; E:K were already 0
$ws.Range("D61").Value = 0

# Cash Flow Statement rows whose E column (second period) reverts to "NA"
$naRows = 83,89,91,94,100,101,102
foreach ($r in $naRows) {
    $addr = "E" + $r
    $ws.Range($addr).Value = "NA"
}

# Row 96 (Dividends Paid) E column reverts to 0 instead of -1700
$ws.Range("E96").Value = 0
